$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cells
$ws.Range("A13").Value = 112116035
$ws.Range("B13").Value = 90669
$ws.Range("E13").Value = 6003297
$ws.Range("Q13").Value = 653795
$ws.Range("R13").Value = 6577004
$ws.Range("S13").Value = 10

# Text cells
$ws.Range("C13").Value = "Ovaliderad"
$ws.Range("D13").Value = "VU"
$ws.Range("F13").Value = "Spricktaggsvamp"
$ws.Range("G13").Value = "Hydnellum glaucopus"
$ws.Range("H13").Value = "(Maas Geest. & Nannf.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("K13").Value = "teleomorf"
$ws.Range("P13").Value = "Svarvartorp ca 400 m SO om, Upl"
$ws.Range("T13").Value = "Stockholm"
$ws.Range("U13").Value = "Ekerö"
$ws.Range("V13").Value = "Uppland"
$ws.Range("W13").Value = "Ekerö"
$ws.Range("AC13").Value = "Kött ganska sprött (ej korkartat), ljusbrunt."
$ws.Range("AI13").Value = "Gles barrskog på sand (både tall och gran)"
$ws.Range("AW13").Value = "Jan Yngve Andersson"
$ws.Range("AX13").Value = "Jan Yngve Andersson"

# Empty-string (inline string) cells
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("AF13").Value = ""
$ws.Range("AT13").Value = ""
$ws.Range("AY13").Value = ""

# Date-like text cells - must stay text, not auto-convert to a date serial.
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2023-09-14"
$ws.Range("Y13").Style = "Normal"

$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2023-09-14"
$ws.Range("AA13").Style = "Normal"

# Boolean cells
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
